$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update column F (想去人数 / "want to go" count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 518
$ws1.Range("F6").Value = 2821
$ws1.Range("F10").Value = 1586
$ws1.Range("F11").Value = 564
$ws1.Range("F24").Value = 241
$ws1.Range("F25").Value = 32
$ws1.Range("F27").Value = 1798
$ws1.Range("F28").Value = 43
$ws1.Range("F30").Value = 97
$ws1.Range("F31").Value = 571

# Sheet "全部类型" (all types) - same updates, rows offset by +1
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 518
$ws4.Range("F7").Value = 2821
$ws4.Range("F11").Value = 1586
$ws4.Range("F12").Value = 564
$ws4.Range("F25").Value = 241
$ws4.Range("F26").Value = 32
$ws4.Range("F28").Value = 1798
$ws4.Range("F29").Value = 43
$ws4.Range("F31").Value = 97
$ws4.Range("F32").Value = 571
